# Regenerate Report for Handback: update the timestamp values recorded
# for the latest handoff/handback xliff generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-30 09:36:23"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-30 09:36:18"
$wsZhCn.Range("K2").Value = "2016-08-30 09:36:36"

# de-de sheet: Correspond Handoff Datetime (shared with Overview value) /
# Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-30 09:36:23"
$wsDeDe.Range("K2").Value = "2016-08-30 09:36:44"
